$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 476, shifting existing row 476..563 down to 477..564.
$ws.Rows.Item(476).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(476, 1).Value = 10
$ws.Cells.Item(476, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(476, 3).Value = "La Araucanía"
$ws.Cells.Item(476, 4).Value = 45015
$ws.Cells.Item(476, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(476, 5).Value = 9
$ws.Cells.Item(476, 6).Value = 100114014
$ws.Cells.Item(476, 7).Value = "Betarraga"
$ws.Cells.Item(476, 8).Value = "Sin especificar"
$ws.Cells.Item(476, 9).Value = "Primera"
$ws.Cells.Item(476, 10).Value = 180
$ws.Cells.Item(476, 11).Value = 9000
$ws.Cells.Item(476, 12).Value = 10000
$ws.Cells.Item(476, 13).Value = 9556
$ws.Cells.Item(476, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(476, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(476, 16).Value = 796
$ws.Cells.Item(476, 17).Value = 12
$ws.Cells.Item(476, 18).Value = "Hortaliza"
